# Fruta / hortaliza, semanal
# Insert two new weekly data rows at the top of the "Mandarina" data block
# (rows 1094-1095), pushing all existing rows in that block down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 1094; this shifts old rows 1094..1169
# down to 1096..1171 (matching the new dimension A1:T1171) while preserving
# their contents and formatting automatically.
$ws.Rows("1094:1095").Insert()

# --- New row 1094 ---
$ws.Range("A1094").Value = 10
$ws.Range("B1094").Value = "Vega Modelo de Temuco"
$ws.Range("C1094").Value = "La Araucanía"
$ws.Range("D1094").Value = 45265
$ws.Range("E1094").Value = 9
$ws.Range("F1094").Value = "Fruta"
$ws.Range("G1094").Value = 100102
$ws.Range("H1094").Value = "Cítricos"
$ws.Range("I1094").Value = 100102004
$ws.Range("J1094").Value = "Mandarina"
$ws.Range("K1094").Value = "Murcott"
$ws.Range("L1094").Value = "Especial"
$ws.Range("M1094").Value = 200
$ws.Range("N1094").Value = 17000
$ws.Range("O1094").Value = 18000
$ws.Range("P1094").Value = 17500
$ws.Range("Q1094").Value = "`$/bandeja 18 kilos"
$ws.Range("R1094").Value = "Región de O'Higgins"
$ws.Range("S1094").Value = 972
$ws.Range("T1094").Value = 18

# --- New row 1095 ---
$ws.Range("A1095").Value = 10
$ws.Range("B1095").Value = "Vega Modelo de Temuco"
$ws.Range("C1095").Value = "La Araucanía"
$ws.Range("D1095").Value = 45265
$ws.Range("E1095").Value = 9
$ws.Range("F1095").Value = "Fruta"
$ws.Range("G1095").Value = 100102
$ws.Range("H1095").Value = "Cítricos"
$ws.Range("I1095").Value = 100102004
$ws.Range("J1095").Value = "Mandarina"
$ws.Range("K1095").Value = "Murcott"
$ws.Range("L1095").Value = "Primera"
$ws.Range("M1095").Value = 2
$ws.Range("N1095").Value = 495000
$ws.Range("O1095").Value = 495000
$ws.Range("P1095").Value = 495000
$ws.Range("Q1095").Value = "`$/bins (450 kilos)"
$ws.Range("R1095").Value = "Región de O'Higgins"
$ws.Range("S1095").Value = 1100
$ws.Range("T1095").Value = 450
